# ---------------------------------------------------------------------------
# mytest.xlsx: the export was re-run against a newer log batch (3 rows instead
# of 1) with an extra pandas index column pair ("Unnamed: 0.1"/"Unnamed: 0") up
# front, a new "pt_time" column, and new trailing timing/t-SNE columns:
# oc_tr_time, centroid_plot_start, total_oc_time, tsne_n_iter, tsne_perplexity.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new leading index columns; this shifts all existing
# columns (and their styling) two slots to the right.
$ws.Range("A:B").Insert()

# Two more data rows are now present (rows 3 and 4); row 2 keeps its place.
# NOTE: cells are written directly rather than via Rows.Insert() so that no
# row formatting gets inherited from row 2 above.

# --- Header row (row 1), columns B..BF ---
$headers = @(
    "Unnamed: 0.1",
    "Unnamed: 0",
    "id",
    "batch_size",
    "padded_seq_len",
    "padded_char_len",
    "logpath",
    "logfilename",
    "pkl_file",
    "tk_file",
    "load_from_pkl",
    "train_ratio",
    "ablation",
    "save_dir",
    "designated_ukc_cls",
    "clean_part_1",
    "clean_part_2",
    "clean_time_1",
    "clean_part_4",
    "clean_time_2",
    "clean_part_6",
    "radius",
    "ocloss",
    "octrf1",
    "ocvalf1",
    "save_padded_num_sequences",
    "char_embedding_size",
    "pt_optimizer",
    "num_classes",
    "pt_loss",
    "manual_color_map",
    "centroid_black",
    "pt_epochs",
    "oc_epochs",
    "ptmodel_name",
    "data_dir",
    "save_ptmodel",
    "pt_wait",
    "ptmodel_path",
    "pt_time",
    "feature_from",
    "f1Known",
    "F1Open",
    "f1_weighted",
    "debug",
    "f1_micro",
    "oc_accu",
    "ukc_label",
    "store_features",
    "oc_wait",
    "oc_lr",
    "oc_optimizer",
    "oc_tr_time",
    "centroid_plot_start",
    "total_oc_time",
    "tsne_n_iter",
    "tsne_perplexity"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# --- Row 2 (columns A..BF) ---
$row2 = @(
    0,
    0,
    1,
    "2022-04-27_06_53_57.955941_b60511fbc5c811ecba05ac8247733f47",
    32,
    32,
    64,
    "C:\ML_data\Logs\BGL.log",
    "BGL.log",
    "C:\Users\Bhujay_ROG\MyDev\OCLog\oclog\BGL\data\bgl_ukc.pkl",
    "C:\Users\Bhujay_ROG\MyDev\OCLog\oclog\BGL\data\bgl_tk.pkl",
    $true,
    0.8,
    500,
    "data",
    5,
    $true,
    $true,
    $true,
    $true,
    $true,
    $true,
    "4.9835577,4.987583,5.033294,4.9719653,4.9852686",
    278.4970092773438,
    0.7935053703064504,
    0.5880684190574047,
    $false,
    50,
    "adam",
    5,
    "categorical_crossentropy",
    $true,
    $true,
    6,
    10,
    "ptmodel",
    "data",
    $true,
    3,
    "data\ptmodel_2022-04-27_06_51_01.210878/",
    40.39199113845825,
    "train_data",
    48.853,
    0,
    0.4859436025389964,
    $false,
    0.4097222222222222,
    40.97,
    7,
    $true,
    3,
    2,
    $null,
    173.7197544574738,
    1651022637.95694,
    177.4426748752594,
    $null,
    $null
)
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}

# --- Row 3 (columns A..BF) ---
$row3 = @(
    1,
    1,
    $null,
    "2022-04-27_07_52_11.033363_d80e62c1c5d011ecba05ac8247733f47",
    32,
    32,
    64,
    "C:\ML_data\Logs\BGL.log",
    "BGL.log",
    "C:\Users\Bhujay_ROG\MyDev\OCLog\oclog\BGL\data\bgl_ukc.pkl",
    "C:\Users\Bhujay_ROG\MyDev\OCLog\oclog\BGL\data\bgl_tk.pkl",
    $true,
    0.8,
    500,
    "data",
    5,
    $true,
    $true,
    $true,
    $true,
    $true,
    $true,
    "4.9835577,4.987583,5.033294,4.9719653,4.9852686",
    278.4970092773438,
    0.7935053703064504,
    0.5880684190574047,
    $false,
    50,
    "adam",
    5,
    "categorical_crossentropy",
    $true,
    $true,
    6,
    10,
    "ptmodel",
    "data",
    $true,
    3,
    "data\ptmodel_2022-04-27_07_49_27.544349/",
    39.82676148414612,
    "train_data",
    48.853,
    0,
    0.4859436025389964,
    $false,
    0.4097222222222222,
    40.97,
    7,
    $true,
    3,
    2,
    $null,
    160.5271420478821,
    1651026131.033363,
    164.7280344963074,
    2000,
    30
)
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

# --- Row 4 (columns A..BF) ---
$row4 = @(
    1,
    $null,
    $null,
    "2022-04-27_08_20_12.307190_c22c69a3c5d411eca8e1ac8247733f47",
    32,
    32,
    64,
    "C:\ML_data\Logs\BGL.log",
    "BGL.log",
    "C:\Users\Bhujay_ROG\MyDev\OCLog\oclog\BGL\data\bgl_ukc.pkl",
    "C:\Users\Bhujay_ROG\MyDev\OCLog\oclog\BGL\data\bgl_tk.pkl",
    $true,
    0.8,
    5000,
    "data",
    5,
    $true,
    $true,
    $true,
    $true,
    $true,
    $true,
    "6.8490844,6.853139",
    1318.860229492188,
    0.9542974079126875,
    0.9979756753212693,
    $false,
    50,
    "adam",
    2,
    "categorical_crossentropy",
    $true,
    $true,
    6,
    50,
    "ptmodel",
    "data",
    $true,
    3,
    "data\ptmodel_2022-04-27_08_06_47.047017/",
    97.81152153015137,
    "train_data",
    66.6667,
    0,
    0.8585069444444444,
    $false,
    0.8585069444444443,
    85.85,
    7,
    $true,
    3,
    2,
    $null,
    794.5880978107452,
    1651027812.30719,
    808.3329162597656,
    2000,
    8
)
for ($i = 0; $i -lt $row4.Length; $i++) {
    $ws.Cells.Item(4, $i + 1).Value = $row4[$i]
}

# ---------------------------------------------------------------------------
# Style touch-up (values only were set above; the bold/bordered "header style"
# needs to be extended onto the cells that did not inherit it automatically):
#  * B1:C1 and BA1:BF1 are header cells outside the range the column-insert
#    shifted forward, so they came in unstyled.
#  * The insert also dragged the old styled index cell (originally A2) along
#    to C2, so C2:C4 need the style cleared and A2:A4 need it (re)applied.
# ---------------------------------------------------------------------------
$ws.Range("D1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$ws.Range("D1").Copy()
$ws.Range("BA1:BF1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)

$ws.Range("D1").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)

# Re-assert the A-column values (PasteSpecial(formats) above does not touch
# values, but keep things explicit/robust to ordering).
$ws.Cells.Item(2, 1).Value = $row2[0]
$ws.Cells.Item(3, 1).Value = $row3[0]
$ws.Cells.Item(4, 1).Value = $row4[0]
